$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value  = 1.28
$ws.Range("H2").Value  = 5.2
$ws.Range("I2").Value  = 9.25
$ws.Range("K2").Value  = 2.62
$ws.Range("L2").Value  = 7.4
$ws.Range("N2").Value  = 9.5
$ws.Range("P2").Value  = 4.85
$ws.Range("T2").Value  = 3.45
$ws.Range("X2").Value  = 7.5
$ws.Range("Z2").Value  = 8.75
$ws.Range("AB2").Value = 22
$ws.Range("AC2").Value = 9.5
$ws.Range("AD2").Value = 10.75
$ws.Range("AI2").Value = 70
$ws.Range("AJ2").Value = 28
$ws.Range("AL2").Value = 100
$ws.Range("AM2").Value = 75
$ws.Range("AN2").Value = 3.3
$ws.Range("AP2").Value = 13.5
$ws.Range("AT2").Value = 3.45
$ws.Range("AW2").Value = 10
$ws.Range("AX2").Value = 50
$ws.Range("AY2").Value = 45
$ws.Range("AZ2").Value = 400
$ws.Range("BA2").Value = 350
$ws.Range("BB2").Value = 500
